# Adds "DS - Binary Search Tree" content to the DSA practice sheet:
#  - restyles/extends rows 208-211 (tail of the "Binary Trees" block)
#  - turns rows 214-235 into the new "Binary Search Trees" block (row numbers,
#    restyled B/C/D, new E "yes" markers, a couple of new F links/notes)
#  - adds new F-column reference links/notes scattered earlier in the sheet
#  - gives a batch of pre-existing F-column URL cells real hyperlinks
#  - moves the sheet selection down to where the new rows were authored

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlLeft = -4131

function Copy-Format($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Rows 208-211: existing "Binary Trees" rows that now also get the "Revise"
# (yes) formatting/marker that the rest of that block already has.
# ---------------------------------------------------------------------------
foreach ($r in 208..211) {
    Copy-Format "B205" "B$r"
    Copy-Format "C205" "C$r"
    Copy-Format "D205" "D$r"
    $ws.Range("D$r").Value2 = "yes"
    $ws.Range("E$r").Value2 = "yes"
}

# ---------------------------------------------------------------------------
# Rows 214-230: the new "Binary Search Trees" topic rows. Each gets a
# question number in column A, the same restyle as above, and a "yes" marker
# in column E.
# ---------------------------------------------------------------------------
$questionNumber = @{
    214 = 156; 215 = 157; 216 = 158; 217 = 159; 218 = 160
    219 = 161; 220 = 162; 221 = 163; 222 = 164; 223 = 165
    224 = 166; 225 = 167; 226 = 168; 227 = 169; 228 = 170
    229 = 171; 230 = 172
}

foreach ($r in 214..230) {
    $ws.Range("A$r").Value2 = $questionNumber[$r]
    Copy-Format "B205" "B$r"
    if ($r -eq 220) {
        # This row already had the left-aligned hyperlink-style variant;
        # keep the left alignment while still picking up the highlight fill.
        Copy-Format "C205" "C$r"
        $ws.Range("C$r").HorizontalAlignment = $xlLeft
    } else {
        Copy-Format "C205" "C$r"
    }
    Copy-Format "D205" "D$r"
    $ws.Range("D$r").Value2 = "yes"
    $ws.Range("E$r").Value2 = "yes"
}

# ---------------------------------------------------------------------------
# Rows 231-235: just gain the running question number in column A; their
# B/C/D formatting is left exactly as-is.
# ---------------------------------------------------------------------------
$tailNumber = @{ 231 = 173; 232 = 174; 233 = 175; 234 = 176; 235 = 177 }
foreach ($r in 231..235) {
    $ws.Range("A$r").Value2 = $tailNumber[$r]
}

# ---------------------------------------------------------------------------
# New F-column content inside the new block.
# ---------------------------------------------------------------------------
$ws.Range("F223").Value2 = "https://www.geeksforgeeks.org/avl-tree-set-1-insertion/"
$ws.Hyperlinks.Add($ws.Range("F223"), $ws.Range("F223").Value2) | Out-Null

$ws.Range("F224").Value2 = "links in comments in code"

$ws.Range("F226").Value2 = "https://www.youtube.com/watch?v=wGXB9OWhPTg&ab_channel=TusharRoy-CodingMadeSimple"
$ws.Hyperlinks.Add($ws.Range("F226"), $ws.Range("F226").Value2) | Out-Null
# This particular link only needs the hyperlink look (style), not an actual
# navigable link, so drop the link object but keep the formatting it applied.
$ws.Range("F226").Hyperlinks.Delete() | Out-Null

# New F-column link further up the sheet (same topic family, row 182).
$ws.Range("F182").Value2 = "https://www.geeksforgeeks.org/inorder-tree-traversal-without-recursion-and-without-stack/"
$ws.Hyperlinks.Add($ws.Range("F182"), $ws.Range("F182").Value2) | Out-Null

# ---------------------------------------------------------------------------
# Existing F-column cells that already held a plain URL/text now become real
# hyperlinks (their own text is the target address).
# ---------------------------------------------------------------------------
$existingLinkCells = @("F16","F21","F39","F40","F51","F161","F164","F165","F187","F303","F306","F316","F322")
foreach ($addr in $existingLinkCells) {
    $ws.Hyperlinks.Add($ws.Range($addr), $ws.Range($addr).Value2) | Out-Null
}

# ---------------------------------------------------------------------------
# Leave the sheet scrolled/selected where the author ended up editing.
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("A226"), $true)
$ws.Range("B230").Select() | Out-Null

Write-Output "done"
